$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# RegressionModel sheet: grows from a 2-row table (A1:C3) to a 9-row table
# (A1:C10), re-sorted alphabetically by floor, with refreshed MAE numbers.
# ---------------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("RegressionModel")

# Extend column-A's numbering style (bold/border/center, same as the header)
# down through the new rows before filling in values.
$wsReg.Range("A2").Copy()
$wsReg.Range("A4:A10").PasteSpecial(-4122)

$regData = @(
    @(0, "Floor 10_aggr", 6.399127039792575),
    @(1, "Floor 12_aggr", 6.713001223166926),
    @(2, "Floor 14_aggr", 7.059946410372695),
    @(3, "Floor 16_aggr", 8.402791164765512),
    @(4, "Floor 3_aggr", 6.950256743968994),
    @(5, "Floor 4_aggr", 6.772114049023483),
    @(6, "Floor 6_aggr", 7.039346536089296),
    @(7, "Floor 8_aggr", 6.98998673392105),
    @(8, "String Pots_aggr", 1.546210047362643)
)

$row = 2
foreach ($r in $regData) {
    $wsReg.Cells.Item($row, 1).Value = $r[0]
    $wsReg.Cells.Item($row, 2).Value = $r[1]
    $wsReg.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# ARIMA sheet: grows from a 1-row table (A1:C2) to a 9-row table (A1:C10),
# same floor ordering, with refreshed MAE numbers.
# ---------------------------------------------------------------------------
$wsArima = $wb.Worksheets.Item("ARIMA")

$wsArima.Range("A2").Copy()
$wsArima.Range("A3:A10").PasteSpecial(-4122)

$arimaData = @(
    @(0, "Floor 10_aggr", 7.289631405653163),
    @(1, "Floor 12_aggr", 6.774843235548796),
    @(2, "Floor 14_aggr", 6.857340676566108),
    @(3, "Floor 16_aggr", 8.51881499337639),
    @(4, "Floor 3_aggr", 6.984120514459121),
    @(5, "Floor 4_aggr", 6.867393864252254),
    @(6, "Floor 6_aggr", 7.098759493396401),
    @(7, "Floor 8_aggr", 6.874483117872179),
    @(8, "String Pots_aggr", 1.184971391881689)
)

$row = 2
foreach ($r in $arimaData) {
    $wsArima.Cells.Item($row, 1).Value = $r[0]
    $wsArima.Cells.Item($row, 2).Value = $r[1]
    $wsArima.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

Write-Host "done"
